# PRIXPS.xlsx - "Add files via upload" edit
#
# 1. D37, D38, D39 all get the corrected unit price of 815.
# 2. The sheet view is scrolled down (row 29 at the top) and the selection
#    moves from A2:D48 to C2:D48 (active cell C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# --- Corrected values in column D ---
$ws.Range("D37").Value = 815
$ws.Range("D38").Value = 815
$ws.Range("D39").Value = 815

# --- View state: scroll + new selection ---
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C2:D48").Select()
